$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order identifiers updated) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961402435634"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961426195998"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961426195998"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961426675618"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961427315965"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961402115612.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961402275846.csv"
$ws1.Range("B4").Value = "go_stims-16509961402275846.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961402435634.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509961426035974.csv"
$ws2.Range("B3").Value = "ZB-match_6-16509961406195607.csv"
$ws2.Range("B4").Value = "OB-1650996141363592.csv"
$ws2.Range("B5").Value = "TB-16509961421235595.csv"
$ws2.Range("B6").Value = "OB-16509961414835615.csv"
$ws2.Range("B7").Value = "OB-1650996141259558.csv"
$ws2.Range("B8").Value = "ZB-match_3-16509961411155593.csv"
$ws2.Range("B9").Value = "ZB-match_9-16509961404675975.csv"
$ws2.Range("B10").Value = "TB-16509961424035614.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961426355963.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961426195998.csv"
$ws4.Range("B4").Value = "MM_stims-16509961426515956.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961426355963.csv"
$ws4.Range("B6").Value = "MM_stims-16509961426675618.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961426515956.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650996142683595.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961426675618.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961426995606.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961427156.csv"
